$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Step 1: drop the bulk of the body ------------------------------------
# Paragraphs 2 ("Group Assignment") through 28 ("Herbert:Safety") are
# removed outright; paragraph 29 ("Cesar:Process") survives (repurposed
# below), as do the two trailing empty paragraphs.
$delStart = $d.Paragraphs(2).Range.Start
$delEnd = $d.Paragraphs(29).Range.Start
$d.Range($delStart, $delEnd).Delete()

# --- Step 2: title paragraph ("Process and Safety" -> merged title) ------
$titleXml = "<w:p $wNs><w:pPr><w:pStyle w:val=`"Title`"/><w:jc w:val=`"center`"/></w:pPr><w:r><w:t>Process and Safety Group Assignment</w:t></w:r></w:p>"
[void]$d.Paragraphs(1).Range.InsertXML($titleXml)

# --- Step 3: "Process" heading (was "Cesar:Process") ---------------------
$processXml = "<w:p $wNs><w:pPr><w:pStyle w:val=`"Heading1`"/></w:pPr><w:r><w:t>Process</w:t></w:r></w:p>"
[void]$d.Paragraphs(2).Range.InsertXML($processXml)

# --- Step 4: two new Heading2 paragraphs ----------------------------------
$funcXml = "<w:p $wNs><w:pPr><w:pStyle w:val=`"Heading2`"/></w:pPr><w:r><w:t>Functionally analyzed process of the project</w:t></w:r></w:p>"
[void]$d.Paragraphs(3).Range.InsertXML($funcXml)

$interXml = "<w:p $wNs><w:pPr><w:pStyle w:val=`"Heading2`"/></w:pPr><w:r><w:t>Interaction between Parts/Modules</w:t></w:r></w:p>"
[void]$d.Paragraphs(4).Range.InsertXML($interXml)

# --- Step 5: two trailing empty paragraphs --------------------------------
$emptyXml = "<w:p $wNs/>"
$endRange = $d.Content
$endRange.Collapse(0)
[void]$endRange.InsertXML($emptyXml)
$endRange = $d.Content
$endRange.Collapse(0)
[void]$endRange.InsertXML($emptyXml)
